$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Species abbreviation codes being filled into column E ("species" id),
# in the order they were originally entered (this affects the order new
# entries land in the shared-string table).
$entries = @(
    @{ Row = 2;  Id = "age_til" },
    @{ Row = 5;  Id = "eur_div" },
    @{ Row = 17; Id = "dep_acr" },
    @{ Row = 21; Id = "ono_sen" },
    @{ Row = 23; Id = "osm_cla" },
    @{ Row = 7;  Id = "lys_cil" },
    @{ Row = 15; Id = "bot_mul" },
    @{ Row = 14; Id = "bot_dis" },
    @{ Row = 18; Id = "dry_cri" },
    @{ Row = 19; Id = "dry_int" },
    @{ Row = 20; Id = "dry_mar" },
    @{ Row = 16; Id = "den_pun" },
    @{ Row = 22; Id = "osm_cin" },
    @{ Row = 28; Id = "den_den" },
    @{ Row = 33; Id = "spi_ann" },
    @{ Row = 24; Id = "phe_con" },
    @{ Row = 25; Id = "pol_acr" },
    @{ Row = 26; Id = "pte_aqu" },
    @{ Row = 27; Id = "the_nov" },
    @{ Row = 29; Id = "den_obs" },
    @{ Row = 30; Id = "dip_com" },
    @{ Row = 31; Id = "dip_dig" },
    @{ Row = 32; Id = "lyc_cla" }
)

# Template cell that already carries the formatting used throughout
# column E (vertically centered, matching font).
$fmtSource = $ws.Range("E1")

foreach ($entry in $entries) {
    $target = $ws.Cells.Item($entry.Row, 5)
    $target.Value = $entry.Id

    $fmtSource.Copy() | Out-Null
    $target.PasteSpecial(-4122) | Out-Null
}

$excel.CutCopyMode = 0

$ws.Range("H24").Select() | Out-Null
